# RecapAsta.xlsx -- refresh team rosters, add section/header styling
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the two team columns
$ws.Columns.Item(1).ColumnWidth = 39.15
$ws.Columns.Item(2).ColumnWidth = 39.15

# Team names + roster rows
$ws.Range("A1").Value = 'Salernitana'
$ws.Range("B1").Value = 'Napoli'
$ws.Range("A2").Value = 'Portieri'
$ws.Range("B2").Value = 'Portieri'
$ws.Range("A3").Value = 'PADELLI  --  1  --  Inter'
$ws.Range("B3").Value = 'RINALDI  --  1  --  Parma'
$ws.Range("A4").Value = 'RADUNOVIC  --  1  --  Atalanta'
$ws.Range("B4").Value = 'SZCZESNY  --  1  --  Juventus'
$ws.Range("A5").Value = 'HANDANOVIC  --  1  --  Inter'
$ws.Range("B5").Value = 'SPORTIELLO  --  1  --  Atalanta'
$ws.Range("A6").Value = 'Difensori'
$ws.Range("B6").Value = 'Difensori'
$ws.Range("A7").Value = 'MARRONE  --  1  --  Crotone'
$ws.Range("B7").Value = 'KUMBULLA  --  1  --  Roma'
$ws.Range("A8").Value = 'GOLDANIGA  --  1  --  Genoa'
$ws.Range("B8").Value = 'IACOPONI  --  1  --  Parma'
$ws.Range("A9").Value = 'RAMOS  --  1  --  Spezia'
$ws.Range("B9").Value = 'TOMORI  --  1  --  Milan'
$ws.Range("A10").Value = 'TRIPALDELLI  --  1  --  Cagliari'
$ws.Range("B10").Value = 'CUADRADO  --  1  --  Juventus'
$ws.Range("A11").Value = 'PEZZELLA GER.  --  1  --  Fiorentina'
$ws.Range("B11").Value = 'YOUNG  --  1  --  Inter'
$ws.Range("A12").Value = 'SMALLING  --  1  --  Roma'
$ws.Range("B12").Value = 'MBAYE  --  1  --  Bologna'
$ws.Range("A13").Value = 'MURRU  --  1  --  Torino'
$ws.Range("B13").Value = 'BASTONI  --  1  --  Inter'
$ws.Range("A14").Value = 'HATEBOER  --  1  --  Atalanta'
$ws.Range("B14").Value = 'RANOCCHIA  --  1  --  Inter'
$ws.Range("A15").Value = 'Centrocampisti'
$ws.Range("B15").Value = 'Centrocampisti'
$ws.Range("A16").Value = 'ANDERSON D.  --  1  --  Lazio'
$ws.Range("B16").Value = 'ERIKSEN  --  1  --  Inter'
$ws.Range("A17").Value = 'BASELLI  --  1  --  Torino'
$ws.Range("B17").Value = 'LOCATELLI  --  1  --  Sassuolo'
$ws.Range("A18").Value = 'SVANBERG  --  1  --  Bologna'
$ws.Range("B18").Value = 'ROJAS  --  1  --  Crotone'
$ws.Range("A19").Value = 'BARELLA  --  1  --  Inter'
$ws.Range("B19").Value = 'PEREIRO  --  1  --  Cagliari'
$ws.Range("A20").Value = 'GRASSI  --  1  --  Parma'
$ws.Range("B20").Value = 'KOVALENKO  --  1  --  Atalanta'
$ws.Range("A21").Value = 'KULUSEVSKI  --  1  --  Juventus'
$ws.Range("B21").Value = 'MELEGONI  --  1  --  Genoa'
$ws.Range("A22").Value = 'CYPRIEN  --  1  --  Parma'
$ws.Range("B22").Value = 'DOMINGUEZ  --  1  --  Bologna'
$ws.Range("A23").Value = 'AGUDELO  --  1  --  Spezia'
$ws.Range("B23").Value = 'MIRANCHUK  --  1  --  Atalanta'
$ws.Range("A24").Value = 'Attaccanti'
$ws.Range("B24").Value = 'Attaccanti'
$ws.Range("A25").Value = 'FARIAS  --  1  --  Spezia'
$ws.Range("B25").Value = 'IMMOBILE  --  1  --  Lazio'
$ws.Range("A26").Value = 'SIMEONE  --  1  --  Cagliari'
$ws.Range("B26").Value = 'KOUAME''  --  1  --  Fiorentina'
$ws.Range("A27").Value = 'RIVIERE  --  1  --  Crotone'
$ws.Range("B27").Value = 'SAU  --  1  --  Benevento'
$ws.Range("A28").Value = 'BRAAF  --  1  --  Udinese'
$ws.Range("B28").Value = 'ZAZA  --  1  --  Torino'
$ws.Range("A29").Value = 'NESTOROVSKI  --  1  --  Udinese'
$ws.Range("B29").Value = 'TORREGROSSA  --  1  --  Sampdoria'
$ws.Range("A30").Value = 'OKAKA  --  1  --  Udinese'
$ws.Range("B30").Value = 'PANDEV  --  1  --  Genoa'

# Title row: bold, size 20
$title = $ws.Range("A1:B1")
$title.Font.Bold = $true
$title.Font.Size = 20

# Section headers: bold, size 16, distinct colors per section
$portieri = $ws.Range("A2:B2")
$portieri.Font.Bold = $true
$portieri.Font.Size = 16
$portieri.Font.Color = 26367      # FFFF6600 orange

$difensori = $ws.Range("A6:B6")
$difensori.Font.Bold = $true
$difensori.Font.Size = 16
$difensori.Font.Color = 16711680  # FF0000FF blue

$centrocampisti = $ws.Range("A15:B15")
$centrocampisti.Font.Bold = $true
$centrocampisti.Font.Size = 16
$centrocampisti.Font.Color = 32768  # FF008000 green

$attaccanti = $ws.Range("A24:B24")
$attaccanti.Font.Bold = $true
$attaccanti.Font.Size = 16
$attaccanti.Font.Color = 255       # FFFF0000 red

Write-Output "RecapAsta updated"
